$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to existing rows 393-425 (columns F and/or G)
$updates = @{
    393 = @{ F = 308141; G = 1242 }
    394 = @{ F = 166243 }
    395 = @{ F = 751141; G = 1957 }
    396 = @{ F = 164905 }
    397 = @{ F = 108147 }
    398 = @{ F = 298509; G = 1470 }
    399 = @{ F = 200644; G = 968 }
    400 = @{ F = 150917; G = 766 }
    401 = @{ F = 273672; G = 935 }
    402 = @{ F = 717510; G = 1389 }
    403 = @{ F = 352311 }
    404 = @{ F = 225161; G = 913 }
    405 = @{ F = 174044; G = 693 }
    406 = @{ F = 171049; G = 682 }
    407 = @{ F = 158120; G = 675 }
    408 = @{ F = 303975; G = 835 }
    409 = @{ F = 703517 }
    410 = @{ F = 363878 }
    411 = @{ F = 225251; G = 828 }
    412 = @{ F = 176069; G = 646 }
    413 = @{ F = 149091; G = 659 }
    414 = @{ F = 146284; G = 554 }
    415 = @{ F = 306610 }
    416 = @{ F = 659038; G = 923 }
    417 = @{ F = 332629; G = 577 }
    418 = @{ F = 200845; G = 699 }
    419 = @{ F = 147632; G = 504 }
    420 = @{ F = 136615; G = 492 }
    421 = @{ F = 150536 }
    422 = @{ F = 291227 }
    423 = @{ F = 426141 }
    424 = @{ F = 250419 }
    425 = @{ F = 136103; G = 536 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# Append new row 426 with data for 2021-05-05
$ws.Range("A426").Value = 44320
$ws.Range("B426").Value = 384317
$ws.Range("C426").Value = 8020
$ws.Range("D426").Value = 708
$ws.Range("E426").Value = 11886
$ws.Range("F426").Value = 90229
$ws.Range("G426").Value = 501

# Match the date formatting style used in column A (numFmt yyyy-mm-dd)
$ws.Range("A426").NumberFormat = "yyyy-mm-dd"
